$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value = "Arica y Parinacota"
$ws.Range("D3").Value = 44428
$ws.Range("E3").Value = 15
$ws.Range("F3").Value = 100112045
$ws.Range("G3").Value = "Zapallo"
$ws.Range("H3").Value = "Camote"
$ws.Range("I3").Value = "1a nueva(o)"
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 580
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = 590
$ws.Range("N3").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O3").Value = "Perú"
$ws.Range("P3").Value = 590
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
